$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*nValues*") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the target paragraph containing 'nValues'"
}

$r = $target.Range
$full = $r.WordOpenXML

$newPara = '<w:p w14:paraId="60CD359E" w14:textId="04DFE644" w:rsidR="00821D22" w:rsidRPr="00821D22" w:rsidRDefault="00821D22" w:rsidP="00821D22"><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:ind w:left="1080"/><w:rPr><w:color w:val="BCBEC4"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="CF8E6D"/></w:rPr><w:t>int</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t>] n</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t>Values</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve"> = {</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>10</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>50</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>75</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>100</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>200</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>225</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>400</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>450</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>500</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>600</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>750</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>800</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>1200</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>1600</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>1800</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>2200</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>2400</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>4000</w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>6000</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:color w:val="2AACB8"/></w:rPr><w:t>8000</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00821D22"><w:rPr><w:color w:val="BCBEC4"/></w:rPr><w:t>};</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

$bodyTag = '<w:body>'
$bodyIdx = $full.IndexOf($bodyTag)
if ($bodyIdx -lt 0) { throw "could not find <w:body> in WordOpenXML" }

$pStartIdx = $full.IndexOf('<w:p ', $bodyIdx)
if ($pStartIdx -lt 0) { throw "could not find opening <w:p in WordOpenXML body" }

$pEndTag = '</w:p>'
$pEndIdx = $full.IndexOf($pEndTag, $pStartIdx)
if ($pEndIdx -lt 0) { throw "could not find closing </w:p> in WordOpenXML body" }
$pEndIdx = $pEndIdx + $pEndTag.Length

$before = $full.Substring(0, $pStartIdx)
$after = $full.Substring($pEndIdx)

$full = $before + $newPara + $after

$r.InsertXML($full)

Write-Output "done"
